$d = $word.ActiveDocument

$replacements = @(
    @("346÷8=", "993÷2="),
    @("293÷9=", "296÷7="),
    @("347÷9=", "310÷8="),
    @("417÷4=", "863÷6="),
    @("266÷9=", "390÷8="),
    @("600÷9=", "805÷3="),
    @("157÷9=", "295÷6="),
    @("177÷2=", "176÷3="),
    @("528÷3=", "672÷3="),
    @("986÷3=", "427÷5="),
    @("884÷2=", "725÷6="),
    @("870÷5=", "519÷9="),
    @("401÷7=", "789÷9="),
    @("598÷7=", "737÷5="),
    @("758÷7=", "550÷7="),
    @("841÷5=", "948÷6="),
    @("453÷4=", "920÷9="),
    @("612÷7=", "930÷6="),
    @("446÷4=", "576÷2="),
    @("128÷2=", "524÷8="),
    @("845÷4=", "684÷4="),
    @("666÷5=", "110÷6="),
    @("701÷7=", "854÷4="),
    @("781÷3=", "451÷3="),
    @("696÷4=", "158÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
